$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 - Average Fitness/Efficiency +/- stDev (text, shared strings)
$ws.Range("C10").Value = "7674.6 ± 100.055"
$ws.Range("D10").Value = "8681.53 ± 24.4254"
$ws.Range("E10").Value = "6476.07 ± 271.445"
$ws.Range("F10").Value = "9410.37 ± 96.6692"
$ws.Range("G10").Value = "92368.5 ± 1734.06"
$ws.Range("H10").Value = "2084.93 ± 26.2651"
$ws.Range("I10").Value = "728.3 ± 53.8633"

# Row 11 - Best case numbers
$ws.Range("D11").Value = 8722
$ws.Range("E11").Value = 6339

# Row 12 - Worst case numbers
$ws.Range("D12").Value = 8632
$ws.Range("E12").Value = 6912
$ws.Range("F12").Value = 9177
$ws.Range("G12").Value = 90444
$ws.Range("H12").Value = 2025
$ws.Range("I12").Value = 603

# Row 13 - Reliability percentages
$ws.Range("D13").Value = 0.033333300000000003
$ws.Range("E13").Value = 0.26666699999999999
$ws.Range("G13").Value = 0.2
$ws.Range("I13").Value = 0.2

# Row 14 - Avg Efficiency/Fitness +/- stDev (text, shared strings)
$ws.Range("C14").Value = "10412.5 ± 4936.15"
$ws.Range("D14").Value = "16000 ± nan"
$ws.Range("E14").Value = "13125 ± 5616.88"
$ws.Range("F14").Value = "16137.5 ± 6844.75"
$ws.Range("G14").Value = "5825 ± 4504.3"
$ws.Range("H14").Value = "17700 ± 70.7107"
$ws.Range("I14").Value = "7425 ± 3602.33"

# Row 15 - Best Eff numbers
$ws.Range("D15").Value = 16000
$ws.Range("E15").Value = 7800
$ws.Range("F15").Value = 8750
$ws.Range("G15").Value = 3150
$ws.Range("H15").Value = 17650
$ws.Range("I15").Value = 4250

# Row 16 - Worse Eff numbers
$ws.Range("D16").Value = 16000
$ws.Range("E16").Value = 21750
$ws.Range("F16").Value = 23950
$ws.Range("G16").Value = 14900
$ws.Range("H16").Value = 17750
$ws.Range("I16").Value = 14150

# Column G width change (target stored width 17.33203125; nearest reachable
# value given column-width pixel quantization in this environment)
$ws.Columns("G").ColumnWidth = 16.5

# Selection change
$ws.Range("H19").Select()
